$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("H5").Value = 324.5
$ws.Range("I5").Value = 170.85715
$ws.Range("J5").Value = 1400
$ws.Range("K5").Value = 170.85715
$ws.Range("L5").Value = 1400
$ws.Range("M5").Value = -55.85714999999999
$ws.Range("N5").Value = -1630
$ws.Range("H17").Value = 1695.4445
$ws.Range("J17").Value = 1695.4445
$ws.Range("L17").Value = 5086.333500000001
$ws.Range("N17").Value = -5422.333500000001
$ws.Range("H29").Value = 4833.6665
$ws.Range("J29").Value = 6545.091
$ws.Range("L29").Value = 19635.273
$ws.Range("N29").Value = -20197.273
$ws.Range("H33").Value = 178.57143
$ws.Range("I33").Value = 180.07692
$ws.Range("K33").Value = 180.07692
$ws.Range("M33").Value = 48.92308
$ws.Range("H38").Value = 137
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H58").Value = 439.83334
$ws.Range("I58").Value = 609.75
$ws.Range("K58").Value = 1829.25
$ws.Range("M58").Value = -1679.25
$ws.Range("H74").Value = 5803.853
$ws.Range("I74").Value = 5139.1333
$ws.Range("J74").Value = 6328.6313
$ws.Range("K74").Value = 5139.1333
$ws.Range("L74").Value = 6328.6313
$ws.Range("M74").Value = -4203.1333
$ws.Range("N74").Value = -8200.631300000001
$ws.Range("H77").Value = 5803.853
$ws.Range("I77").Value = 5139.1333
$ws.Range("J77").Value = 6328.6313
$ws.Range("K77").Value = 25695.6665
$ws.Range("L77").Value = 31643.1565
$ws.Range("M77").Value = -21015.6665
$ws.Range("N77").Value = -41003.1565
$ws.Range("H129").Value = 1940.5714
$ws.Range("I129").Value = 1197.6666
$ws.Range("K129").Value = 3592.9998
$ws.Range("M129").Value = 1407.0002

$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 3151.1936
$ws.Range("I32").Value = 2294.9153
$ws.Range("K32").Value = 2294.9153
$ws.Range("M32").Value = -2007.9153
$ws.Range("H74").Value = 2601.0938
$ws.Range("I74").Value = 2512.4644
$ws.Range("K74").Value = 2512.4644
$ws.Range("M74").Value = -1638.4644
$ws.Range("H77").Value = 2601.0938
$ws.Range("I77").Value = 2512.4644
$ws.Range("K77").Value = 12562.322
$ws.Range("M77").Value = -8194.322
$ws.Range("H88").Value = 11535.308
$ws.Range("I88").Value = 1898.3334
$ws.Range("J88").Value = 14426.4
$ws.Range("K88").Value = 1898.3334
$ws.Range("L88").Value = 14426.4
$ws.Range("M88").Value = -1492.3334
$ws.Range("N88").Value = -15238.4
$ws.Range("H91").Value = 11535.308
$ws.Range("I91").Value = 1898.3334
$ws.Range("J91").Value = 14426.4
$ws.Range("K91").Value = 1898.3334
$ws.Range("L91").Value = 14426.4
$ws.Range("M91").Value = -494.3334
$ws.Range("N91").Value = -17234.4
$ws.Range("H110").Value = 413.8
$ws.Range("I110").Value = 413.8
$ws.Range("K110").Value = 413.8
$ws.Range("M110").Value = 1631.2
$ws.Range("H132").Value = 3537.9534
$ws.Range("I132").Value = 3390.85
$ws.Range("K132").Value = 10172.55
$ws.Range("M132").Value = -7642.549999999999

$ws = $wb.Worksheets.Item(3)
$ws.Range("H105").Value = 3008.4546
$ws.Range("I105").Value = 2313.9412
$ws.Range("K105").Value = 2313.9412
$ws.Range("M105").Value = -566.9412000000002
$ws.Range("H107").Value = 1659.7
$ws.Range("J107").Value = 2546.5
$ws.Range("L107").Value = 2546.5
$ws.Range("N107").Value = -6386.5

$ws = $wb.Worksheets.Item(4)
$ws.Range("H41").Value = 21424.334
$ws.Range("J41").Value = 49990
$ws.Range("L41").Value = 49990
$ws.Range("N41").Value = -50846
$ws.Range("H58").Value = 2116.125
$ws.Range("I58").Value = 1342.8462
$ws.Range("K58").Value = 1342.8462
$ws.Range("M58").Value = -1139.8462
$ws.Range("H136").Value = 2116.125
$ws.Range("I136").Value = 1342.8462
$ws.Range("K136").Value = 4028.5386
$ws.Range("M136").Value = -1478.5386
$ws.Range("H141").Value = 265322.6
$ws.Range("J141").Value = 285914
$ws.Range("L141").Value = 285914
$ws.Range("N141").Value = -296274

$ws = $wb.Worksheets.Item(5)
$ws.Range("H68").Value = 2708
$ws.Range("I68").Value = 2275
$ws.Range("J68").Value = 2996.6667
$ws.Range("K68").Value = 6825
$ws.Range("L68").Value = 8990.000100000001
$ws.Range("M68").Value = -6014
$ws.Range("N68").Value = -10612.0001
$ws.Range("H71").Value = 2708
$ws.Range("I71").Value = 2275
$ws.Range("J71").Value = 2996.6667
$ws.Range("K71").Value = 20475
$ws.Range("L71").Value = 26970.0003
$ws.Range("M71").Value = -16419
$ws.Range("N71").Value = -35082.0003
$ws.Range("H86").Value = 5180
$ws.Range("J86").Value = 16916.5
$ws.Range("L86").Value = 50749.5
$ws.Range("N86").Value = -53121.5
$ws.Range("H89").Value = 5180
$ws.Range("J89").Value = 16916.5
$ws.Range("L89").Value = 152248.5
$ws.Range("N89").Value = -164104.5
$ws.Range("H92").Value = 308.5
$ws.Range("J92").Value = 586
$ws.Range("L92").Value = 1758
$ws.Range("N92").Value = -4254
$ws.Range("H141").Value = 33242.133
$ws.Range("I141").Value = 9283
$ws.Range("K141").Value = 27849
$ws.Range("M141").Value = -22669

$ws = $wb.Worksheets.Item(6)
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H24").Value = 11749.73
$ws.Range("J24").Value = 9791.25
$ws.Range("L24").Value = 9791.25
$ws.Range("N24").Value = -10137.25
$ws.Range("H29").Value = 6833
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 6833
$ws.Range("K29").Value = 0
$ws.Range("M29").Value = 6833
$ws.Range("N29").Value = -7413
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("N62").Value = 0
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("N65").Value = 0
$ws.Range("H80").Value = 4299.5
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4299.5
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = 4299.5
$ws.Range("N80").Value = -6295.5
$ws.Range("H83").Value = 4299.5
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4299.5
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = 21497.5
$ws.Range("N83").Value = -31481.5
$ws.Range("H97").Value = 956.21875
$ws.Range("I97").Value = 643.6087
$ws.Range("K97").Value = 643.6087
$ws.Range("M97").Value = -147.6087
$ws.Range("H126").Value = 7228.4
$ws.Range("I126").Value = 7822.4165
$ws.Range("J126").Value = 6337.375
$ws.Range("K126").Value = 23467.2495
$ws.Range("L126").Value = 19012.125
$ws.Range("M126").Value = -20997.2495
$ws.Range("N126").Value = -23952.125
$ws.Range("H132").Value = 1741.9556
$ws.Range("I132").Value = 1033.1714
$ws.Range("K132").Value = 3099.5142
$ws.Range("M132").Value = -569.5141999999996
$ws.Range("H134").Value = 66296.336
$ws.Range("J134").Value = 66296.336
$ws.Range("L134").Value = 198889.008
$ws.Range("N134").Value = -203959.008
$ws.Range("H136").Value = 24122.762
$ws.Range("J136").Value = 24122.762
$ws.Range("L136").Value = 72368.28599999999
$ws.Range("N136").Value = -77468.28599999999

$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 1245.2222
$ws.Range("I22").Value = 1011.375
$ws.Range("J22").Value = 1432.3
$ws.Range("K22").Value = 1011.375
$ws.Range("L22").Value = 1432.3
$ws.Range("M22").Value = -716.375
$ws.Range("N22").Value = -2022.3
$ws.Range("H27").Value = 1245.2222
$ws.Range("I27").Value = 1011.375
$ws.Range("J27").Value = 1432.3
$ws.Range("K27").Value = 1011.375
$ws.Range("L27").Value = 1432.3
$ws.Range("M27").Value = -904.375
$ws.Range("N27").Value = -1646.3
$ws.Range("H135").Value = 105056.8
$ws.Range("J135").Value = 105618.664
$ws.Range("L135").Value = 105618.664
$ws.Range("N135").Value = -115758.664

$ws = $wb.Worksheets.Item(8)
$ws.Range("H75").Value = 100000
$ws.Range("J75").Value = 100000
$ws.Range("L75").Value = 100000
$ws.Range("N75").Value = -101872
$ws.Range("H78").Value = 100000
$ws.Range("J78").Value = 100000
$ws.Range("L78").Value = 300000
$ws.Range("N78").Value = -309360
$ws.Range("H94").Value = 45748
$ws.Range("J94").Value = 45748
$ws.Range("L94").Value = 45748
$ws.Range("N94").Value = -47550
$ws.Range("H132").Value = 1512.6471
$ws.Range("I132").Value = 1408.2858
$ws.Range("K132").Value = 4224.857400000001
$ws.Range("M132").Value = -1694.857400000001
